$wb = $excel.ActiveWorkbook
$wsHeader = $wb.Worksheets.Item("Header")
$wsCost = $wb.Worksheets.Item("CostTypes")

# Update the CostTypes sheet data: B5 becomes text "X", C5-C8 get new values
$wsCost.Range("B5").Value = "X"
$wsCost.Range("C5").Value = "Material"
$wsCost.Range("C6").Value = "Personal"
$wsCost.Range("C7").Value = "Test"
$wsCost.Range("C8").Value = "Demo"

# Widen column B on CostTypes sheet
$wsCost.Columns.Item(2).ColumnWidth = 20.44140625

# Update selection on CostTypes sheet and scroll position
$wsCost.Application.Goto($wsCost.Range("A4"))
$wsCost.Range("B5").Select()

# Activate CostTypes sheet (making it the active/selected tab)
$wsCost.Activate()
